$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H40").Value = 4473.913
$ws_ALC.Range("I40").Value = 4999.6665
$ws_ALC.Range("J40").Value = 4288.353
$ws_ALC.Range("K40").Value = 4999.6665
$ws_ALC.Range("L40").Value = 4288.353
$ws_ALC.Range("M40").Value = -4824.6665
$ws_ALC.Range("N40").Value = -4638.353
$ws_ALC.Range("H76").Value = 4453.978
$ws_ALC.Range("I76").Value = 3580.0967
$ws_ALC.Range("J76").Value = 6260
$ws_ALC.Range("K76").Value = 3580.0967
$ws_ALC.Range("L76").Value = 6260
$ws_ALC.Range("M76").Value = -3265.0967
$ws_ALC.Range("N76").Value = -6890
$ws_ALC.Range("H79").Value = 4453.978
$ws_ALC.Range("I79").Value = 3580.0967
$ws_ALC.Range("J79").Value = 6260
$ws_ALC.Range("K79").Value = 3580.0967
$ws_ALC.Range("L79").Value = 6260
$ws_ALC.Range("M79").Value = -2488.0967
$ws_ALC.Range("N79").Value = -8444
$ws_ALC.Range("H111").Value = 2000.4762
$ws_ALC.Range("I111").Value = 3072.25
$ws_ALC.Range("J111").Value = 1340.9231
$ws_ALC.Range("K111").Value = 9216.75
$ws_ALC.Range("L111").Value = 4022.7693
$ws_ALC.Range("M111").Value = -6149.75
$ws_ALC.Range("N111").Value = -10156.7693
$ws_ALC.Range("H125").Value = 1928.2858
$ws_ALC.Range("J125").Value = 706
$ws_ALC.Range("L125").Value = 6354
$ws_ALC.Range("N125").Value = -11274
$ws_ALC.Range("H132").Value = 6417105
$ws_ALC.Range("J132").Value = 25010708
$ws_ALC.Range("L132").Value = 75032124
$ws_ALC.Range("N132").Value = -75037184
$ws_ALC.Range("H138").Value = 10419355
$ws_ALC.Range("I138").Value = 2238.2
$ws_ALC.Range("J138").Value = 27781216
$ws_ALC.Range("K138").Value = 6714.599999999999
$ws_ALC.Range("L138").Value = 83343648
$ws_ALC.Range("M138").Value = -1574.599999999999
$ws_ALC.Range("N138").Value = -83353928
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H5").Value = 980.4
$ws_ARM.Range("I5").Value = 450
$ws_ARM.Range("K5").Value = 450
$ws_ARM.Range("M5").Value = -338
$ws_ARM.Range("H113").Value = 47000
$ws_ARM.Range("J113").Value = 47000
$ws_ARM.Range("L113").Value = 47000
$ws_ARM.Range("N113").Value = -55678
$ws_ARM.Range("H134").Value = 49751.4
$ws_ARM.Range("J134").Value = 49751.4
$ws_ARM.Range("L134").Value = 49751.4
$ws_ARM.Range("N134").Value = -59891.4
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H4").Value = 980.4
$ws_BSM.Range("I4").Value = 450
$ws_BSM.Range("K4").Value = 450
$ws_BSM.Range("M4").Value = -335
$ws_BSM.Range("H20").Value = 990.06665
$ws_BSM.Range("I20").Value = 1215.7778
$ws_BSM.Range("J20").Value = 651.5
$ws_BSM.Range("K20").Value = 1215.7778
$ws_BSM.Range("L20").Value = 651.5
$ws_BSM.Range("M20").Value = -968.7778000000001
$ws_BSM.Range("N20").Value = -1145.5
$ws_BSM.Range("H36").Value = 1893.5
$ws_BSM.Range("I36").Value = 1893.5
$ws_BSM.Range("K36").Value = 1893.5
$ws_BSM.Range("M36").Value = -1359.5
$ws_BSM.Range("H44").Value = 0
$ws_BSM.Range("J44").Value = 0
$ws_BSM.Range("L44").Value = 0
$ws_BSM.Range("N44").ClearContents()
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H58").Value = 3419.1177
$ws_CRP.Range("I58").Value = 2021.2727
$ws_CRP.Range("J58").Value = 5981.8335
$ws_CRP.Range("K58").Value = 2021.2727
$ws_CRP.Range("L58").Value = 5981.8335
$ws_CRP.Range("M58").Value = -1818.2727
$ws_CRP.Range("N58").Value = -6387.8335
$ws_CRP.Range("H134").Value = 702652.9399999999
$ws_CRP.Range("I134").Value = 2707.913
$ws_CRP.Range("K134").Value = 8123.739
$ws_CRP.Range("M134").Value = -5588.739
$ws_CRP.Range("H136").Value = 3419.1177
$ws_CRP.Range("I136").Value = 2021.2727
$ws_CRP.Range("J136").Value = 5981.8335
$ws_CRP.Range("K136").Value = 6063.8181
$ws_CRP.Range("L136").Value = 17945.5005
$ws_CRP.Range("M136").Value = -3513.8181
$ws_CRP.Range("N136").Value = -23045.5005
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H68").Value = 1206.4667
$ws_CUL.Range("J68").Value = 1399.28
$ws_CUL.Range("L68").Value = 4197.84
$ws_CUL.Range("N68").Value = -5819.84
$ws_CUL.Range("H71").Value = 1206.4667
$ws_CUL.Range("J71").Value = 1399.28
$ws_CUL.Range("L71").Value = 12593.52
$ws_CUL.Range("N71").Value = -20705.52
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 28171.72
$ws_GSM.Range("I70").Value = 103446.336
$ws_GSM.Range("J70").Value = 4400.7896
$ws_GSM.Range("K70").Value = 103446.336
$ws_GSM.Range("L70").Value = 4400.7896
$ws_GSM.Range("M70").Value = -103176.336
$ws_GSM.Range("N70").Value = -4940.7896
$ws_GSM.Range("H73").Value = 28171.72
$ws_GSM.Range("I73").Value = 103446.336
$ws_GSM.Range("J73").Value = 4400.7896
$ws_GSM.Range("K73").Value = 103446.336
$ws_GSM.Range("L73").Value = 4400.7896
$ws_GSM.Range("M73").Value = -102510.336
$ws_GSM.Range("N73").Value = -6272.7896
$ws_GSM.Range("H80").Value = 13182903
$ws_GSM.Range("I80").Value = 33335984
$ws_GSM.Range("J80").Value = 1986746.2
$ws_GSM.Range("K80").Value = 33335984
$ws_GSM.Range("L80").Value = 1986746.2
$ws_GSM.Range("M80").Value = -33334986
$ws_GSM.Range("N80").Value = -1988742.2
$ws_GSM.Range("H83").Value = 13182903
$ws_GSM.Range("I83").Value = 33335984
$ws_GSM.Range("J83").Value = 1986746.2
$ws_GSM.Range("K83").Value = 166679920
$ws_GSM.Range("L83").Value = 9933731
$ws_GSM.Range("M83").Value = -166674928
$ws_GSM.Range("N83").Value = -9943715
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 12100
$ws_LTW.Range("I40").Value = 18266.666
$ws_LTW.Range("J40").Value = 7475
$ws_LTW.Range("K40").Value = 18266.666
$ws_LTW.Range("L40").Value = 7475
$ws_LTW.Range("M40").Value = -18130.666
$ws_LTW.Range("N40").Value = -7747
$ws_LTW.Range("H61").Value = 1570.1111
$ws_LTW.Range("I61").Value = 1533
$ws_LTW.Range("J61").Value = 1700
$ws_LTW.Range("K61").Value = 1533
$ws_LTW.Range("L61").Value = 1700
$ws_LTW.Range("M61").Value = -1331
$ws_LTW.Range("N61").Value = -2104
$ws_LTW.Range("H113").Value = 1570.1111
$ws_LTW.Range("I113").Value = 1533
$ws_LTW.Range("J113").Value = 1700
$ws_LTW.Range("K113").Value = 1533
$ws_LTW.Range("L113").Value = 1700
$ws_LTW.Range("M113").Value = 637
$ws_LTW.Range("N113").Value = -6040
$ws_LTW.Range("H140").Value = 45535.332
$ws_LTW.Range("J140").Value = 45535.332
$ws_LTW.Range("L140").Value = 45535.332
$ws_LTW.Range("N140").Value = -55895.332
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H30").Value = 1000
$ws_WVR.Range("I30").Value = 0
$ws_WVR.Range("J30").Value = 1000
$ws_WVR.Range("K30").Value = 0
$ws_WVR.Range("L30").Value = 1000
$ws_WVR.Range("M30").ClearContents()
$ws_WVR.Range("N30").Value = -1214
$ws_WVR.Range("H62").Value = 11207.429
$ws_WVR.Range("I62").Value = 5300.6665
$ws_WVR.Range("J62").Value = 12818.363
$ws_WVR.Range("K62").Value = 5300.6665
$ws_WVR.Range("L62").Value = 12818.363
$ws_WVR.Range("M62").Value = -4676.6665
$ws_WVR.Range("N62").Value = -14066.363
$ws_WVR.Range("H65").Value = 11207.429
$ws_WVR.Range("I65").Value = 5300.6665
$ws_WVR.Range("J65").Value = 12818.363
$ws_WVR.Range("K65").Value = 26503.3325
$ws_WVR.Range("L65").Value = 64091.815
$ws_WVR.Range("M65").Value = -23383.3325
$ws_WVR.Range("N65").Value = -70331.815
$ws_WVR.Range("H126").Value = 4004.2
$ws_WVR.Range("I126").Value = 2226.889
$ws_WVR.Range("K126").Value = 6680.667
$ws_WVR.Range("M126").Value = -4210.667
